$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 15: new diary entry for 1/17/2020 ---
# Column A already holds the literal text "1/17/2020" on row 14 (shared string).
# Assigning that text directly via .Value would be auto-converted to a real
# Excel date serial, so instead copy the already-typed text cell over.
$ws.Range("A14").Copy()
$ws.Range("A15").PasteSpecial(-4104)

# Column B holds a time-of-day fraction; copy the existing "Time" style first
# (xlPasteFormats) so the destination keeps the same number format style as
# the other populated rows, then set the actual value.
$ws.Range("B14").Copy()
$ws.Range("B15").PasteSpecial(-4122)
$ws.Range("B15").Value = 0.925694444444444

$ws.Range("C15").Value = "None"
$ws.Range("D15").Value = "Change the team project to a more suitable one"
$ws.Range("E15").Value = "Finished Goal"
$ws.Range("F15").Value = "I just think that last project is not that appropriate. Anyway, I find another Android app since we can run it and see the features, it's easier and more interesting to explore this project."
$ws.Range("G15").Value = "Not bad"

$ws.Rows("15:15").RowHeight = 102

# --- Row 16: new diary entry for 1/18/2020 ---
# "1/18/2020" has never been typed anywhere else in the sheet, so writing it
# straight to a General-formatted cell would again get auto-parsed into a
# date serial. Build it as a text formula first, then collapse the formula
# down to its plain text result (xlPasteValues) which keeps the General
# style but drops the formula, leaving a plain shared-string cell.
$ws.Range("A16").Formula = "=""1/18/2020"""
$ws.Range("A16").Copy()
$ws.Range("A16").PasteSpecial(-4163)

$ws.Range("B14").Copy()
$ws.Range("B16").PasteSpecial(-4122)
$ws.Range("B16").Value = 0.518055555555556

$ws.Range("C16").Value = "None"
$ws.Range("D16").Value = "Finish the assignment"
$ws.Range("E16").Value = "Finished Goal"
$ws.Range("F16").Value = "It's interesting to explore the Pacman project. I did just what we learnt on class, it's really useful."
$ws.Range("G16").Value = "Not bad"

$ws.Rows("16:16").RowHeight = 64

# --- Update the view so the active cell matches the edit ---
$ws.Range("F16").Select()
